$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Remove the header row "grandes regiões e unidades da federação" (row 6),
# which had no data values of its own; this shifts all subsequent rows up by one.
$ws.Rows("6:6").Delete()
